$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Total" column header in X1
$ws.Range("X1").Value = "Total"

# Row-wise totals for existing rows 2-6 (sum of B:W per row)
$ws.Range("X2").Value = 1996
$ws.Range("X3").Value = 279
$ws.Range("X4").Value = 1120
$ws.Range("X5").Value = 281
$ws.Range("X6").Value = 1419

# New row 7: "Outros" category
$ws.Range("A7").Value = "Outros"
$row7 = @(143,3,20,55,63,80,101,95,95,121,113,121,107,113,107,137,156,164,114,36,13,1)
for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, 2 + $i).Value = $row7[$i]
}
$ws.Range("X7").Value = 1958

# New row 8: "Total" category (column sums of rows 2-7)
$ws.Range("A8").Value = "Total"
$row8 = @(159,11,26,65,72,94,136,157,193,278,358,456,548,608,702,814,845,803,524,159,44,1)
for ($i = 0; $i -lt $row8.Length; $i++) {
    $ws.Cells.Item(8, 2 + $i).Value = $row8[$i]
}
$ws.Range("X8").Value = 7053
